$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1011, 479, 'V', 'Indianapolis', 7, 3, 12, 3, 25, 49, 46.5, 290, 24),
    @(1011, 480, 'H', 'Baltimore', 0, 3, 6, 16, 31, 6, 8.5, -350, 5),
    @(1014, 109, 'V', 'TampaBay', 14, 7, 7, 0, 28, 6.5, 7.5, -290, 0.5),
    @(1014, 110, 'H', 'Philadelphia', 7, 0, 7, 8, 22, 52.5, 53, 245, 26),
    @(1017, 251, 'N', 'Miami', 7, 6, 0, 7, 20, 3.5, 'pk', -125, 24),
    @(1017, 252, 'N', 'Jacksonville', 3, 7, 7, 6, 23, 43.5, 47, 105, 0.5),
    @(1017, 253, 'V', 'Houston', 0, 3, 0, 0, 3, 43.5, 44.5, 450, 21),
    @(1017, 254, 'H', 'Indianapolis', 7, 3, 14, 7, 31, 11, 12.5, -600, 5.5),
    @(1017, 255, 'V', 'GreenBay', 0, 10, 7, 7, 24, 4.5, 5.5, -240, 2),
    @(1017, 256, 'H', 'Chicago', 7, 0, 0, 7, 14, 46, 44, 200, 21.5),
    @(1017, 257, 'V', 'KansasCity', 7, 3, 7, 14, 31, 6, 7, -280, 6.5),
    @(1017, 258, 'H', 'Washington', 3, 10, 0, 0, 13, 54, 54, 240, 27),
    @(1017, 259, 'V', 'Minnesota', 6, 6, 13, 3, 34, 47.5, 1, -135, 1),
    @(1017, 260, 'H', 'Carolina', 7, 3, 7, 11, 28, 1, 45, 115, 21.5),
    @(1017, 261, 'V', 'LAChargers', 0, 6, 0, 0, 6, 48, 51, 135, 0.5),
    @(1017, 262, 'H', 'Baltimore', 7, 10, 10, 7, 34, 3.5, 3, -155, 26),
    @(1017, 263, 'V', 'Cincinnati', 7, 3, 10, 14, 34, 3, 3.5, -185, 2.5),
    @(1017, 264, 'H', 'Detroit', 0, 0, 0, 11, 11, 49, 46.5, 165, 21),
    @(1017, 265, 'V', 'LARams', 0, 28, 3, 7, 38, 6.5, 8.5, -350, 4),
    @(1017, 266, 'H', 'NYGiants', 3, 0, 0, 8, 11, 50, 49.5, 290, 22),
    @(1017, 267, 'V', 'Arizona', 7, 16, 7, 7, 37, 53, 48, 120, 25.5),
    @(1017, 268, 'H', 'Cleveland', 0, 14, 0, 0, 14, 3, 2.5, -140, 3.5),
    @(1017, 269, 'V', 'LasVegas', 10, 7, 14, 3, 34, 44, 45, 180, 23),
    @(1017, 270, 'H', 'Denver', 7, 0, 3, 14, 24, 2.5, 5, -210, 4.5),
    @(1017, 271, 'V', 'Dallas', 7, 3, 7, 12, 35, 1, 3.5, -190, 4),
    @(1017, 272, 'H', 'NewEngland', 14, 0, 0, 15, 29, 48, 50.5, 170, 27),
    @(1017, 273, 'V', 'Seattle', 0, 0, 14, 6, 20, 2.5, 43, 200, 20),
    @(1017, 274, 'H', 'Pittsburgh', 0, 14, 3, 3, 23, 48, 5.5, -240, 1.5),
    @(1018, 275, 'V', 'Buffalo', 3, 17, 11, 0, 31, 3.5, 7, -265, 3.5),
    @(1018, 276, 'H', 'Tennessee', 0, 17, 7, 10, 34, 53, 53.5, 225, 26)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 160 + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}